$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.624.36'
$ws.Range('D3').Value = '1.643.67'
$ws.Range('E3').Value = '  +0.65%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').Value = '215.92'
$ws.Range('E5').Value = '  +1.37%  '
$ws.Range('E6').Value = '  +0.81%  '
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('E9').Value = '  +0.69%  '
$ws.Range('D10').Value = '19.22'
$ws.Range('E10').Value = '  +0.29%  '
$ws.Range('E11').Value = '  -0.09%  '
$ws.Range('D12').Value = '1.873.66'
$ws.Range('E12').Value = '  +0.68%  '
$ws.Range('E13').Value = '  +3.37%  '
$ws.Range('D14').Value = '1.640.60'
$ws.Range('E14').Value = '  +0.34%  '
$ws.Range('E15').Value = '  +1.49%  '
$ws.Range('D16').Value = '65.93'
$ws.Range('E16').Value = '  +4.23%  '
$ws.Range('D17').Value = '26.669.65'
$ws.Range('E18').Value = '  +1.42%  '
$ws.Range('D19').Value = '218.49'
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('E20').Value = '  +0.27%  '
$ws.Range('D21').Value = '4.37'
$ws.Range('E22').Value = '  +2.09%  '
$ws.Range('D23').Value = '9.56'
$ws.Range('E23').Value = '  +1.83%  '
$ws.Range('D24').Value = '2.15'
$ws.Range('E24').Value = '  +11.73%  '
$ws.Range('D25').Value = '146.32'
$ws.Range('E25').Value = '  -1.17%  '
$ws.Range('E26').Value = '  +0.24%  '
$ws.Range('E27').Value = '  -0.39%  '
$ws.Range('E28').Value = '  +3.45%  '
$ws.Range('D29').Value = '15.85'
$ws.Range('E29').Value = '  +2.47%  '
$ws.Range('E30').Value = '  +2.74%  '
$ws.Range('E31').Value = '  +0.98%  '
$ws.Range('D32').Value = '3.40'
$ws.Range('E32').Value = '  +3.24%  '
$ws.Range('D33').Value = '3.05'
$ws.Range('E33').Value = '  +2.68%  '
$ws.Range('D34').Value = '1.277.94'
$ws.Range('E34').Value = '  +5.63%  '
$ws.Range('E35').Value = '  +2.48%  '
$ws.Range('E36').Value = '  +6.21%  '
$ws.Range('E37').Value = '  +0.11%  '
$ws.Range('D38').Value = '0.828'
$ws.Range('E38').Value = '  +2.24%  '
$ws.Range('E39').Value = '  +4.98%  '
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.21%  '
$ws.Range('D41').Value = '0.807'
$ws.Range('E41').Value = '  +2.01%  '
$ws.Range('E42').Value = '  -1.32%  '
$ws.Range('E43').Value = '  +0.83%  '
$ws.Range('D44').Value = '1.784.76'
$ws.Range('E44').Value = '  +0.60%  '
$ws.Range('D45').Value = '93.07'
$ws.Range('E45').Value = '  +0.31%  '
$ws.Range('D46').Value = '59.81'
$ws.Range('E46').Value = '  +9.36%  '
$ws.Range('E47').Value = '  +3.59%  '
$ws.Range('E48').Value = '  +0.68%  '
$ws.Range('E49').Value = '  +1.74%  '
$ws.Range('E50').Value = '  +4.07%  '
$ws.Range('E51').Value = '  -0.70%  '
